$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update geo_abbrvs values (column B) and geo_hierarchy values (column C)
# to restore the values from commit 4fab33f

$ws.Range("B6").Value = "co"        # county
$ws.Range("B11").Value = "anrc"     # alaska native regional corporation
$ws.Range("B12").Value = "aiannh"   # american indian area/alaska native area/hawaiian home land
$ws.Range("B19").Value = "elsd"     # school district (elementary)
$ws.Range("B20").Value = "scsd"     # school district (secondary)
$ws.Range("B21").Value = "unsd"     # school district (unified)
$ws.Range("C23").Value = "state"    # zip code tabulation area hierarchy
$ws.Range("B24").Value = "sldu"     # state legislative district (upper chamber)
$ws.Range("B25").Value = "sldl"     # state legislative district (lower chamber)

# Update view/selection state to reflect scrolled position and active cell
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E19").Select()

$wb.Save()
